$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for column G, rows 2-83, replacing old "Strike#" derived values
$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 2
    6 = 1
    7 = 0
    8 = 2
    9 = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 2
    18 = 2
    19 = 2
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 0
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 3
    37 = 0
    38 = 0
    39 = 1
    40 = 2
    41 = 2
    42 = 2
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    47 = 3
    48 = 0
    49 = 2
    50 = 1
    51 = 2
    52 = 0
    53 = 1
    54 = 0
    55 = 2
    56 = 2
    57 = 0
    58 = 1
    59 = 1
    60 = 1
    61 = 2
    62 = 1
    63 = 3
    64 = 2
    65 = 1
    66 = 0
    67 = 0
    68 = 0
    69 = 1
    70 = 0
    71 = 1
    72 = 1
    73 = 1
    74 = 1
    75 = 2
    76 = 3
    77 = 0
    78 = 0
    79 = 2
    80 = 1
    81 = 1
    82 = 1
    83 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Count) K values in column G"